# Add big-trade-events support to run_config.xlsx
#  - DATASETS: add ES_BIG_TRADES (row 14) and ES_BIG_TRADES_PROXY (row 16) rows
#  - INSTRUMENTS: insert 3 new columns (big_trades_dataset_id,
#    big_trades_proxy_dataset_id, big_trades_source_mode) before the
#    existing "notes" column, and populate the ES row with the new ids.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DATASETS sheet — new rows 14 and 16 (row 13 and 15 stay blank spacers,
# matching the blank-row-separated grouping style already used in this
# sheet between related dataset definitions).
# ---------------------------------------------------------------------
$datasets = $wb.Worksheets.Item("DATASETS")

$datasets.Range("A14").Value = "ES_BIG_TRADES"
$datasets.Range("B14").Value = "big_trades"
$datasets.Range("C14").Value = "canonical"
$datasets.Range("D14").Value = "DB_ES_TRADES"
$datasets.Range("E14").Value = "on_the_fly"
$datasets.Range("F14").Value = "ts_event"
$datasets.Range("G14").Value = "UTC"
$datasets.Range("I14").Value = "event_time"
$datasets.Range("J14").Value = 0
$datasets.Range("K14").Value = "big_trade_events"
$datasets.Range("L14").Value = "instrument_id,session,date"
$datasets.Range("M14").Value = "instrument_id: ES`nthreshold_method: fixed_count`nmin_size: 50"

$datasets.Range("A16").Value = "ES_BIG_TRADES_PROXY"
$datasets.Range("B16").Value = "big_trades_proxy"
$datasets.Range("C16").Value = "canonical"
$datasets.Range("D16").Value = "DB_ES_OHLCV_1S"
$datasets.Range("E16").Value = "on_the_fly"
$datasets.Range("F16").Value = "ts_event"
$datasets.Range("G16").Value = "UTC"
$datasets.Range("H16").Value = "1s"
$datasets.Range("I16").Value = "event_time"
$datasets.Range("J16").Value = 0
$datasets.Range("K16").Value = "big_trade_events_proxy"
$datasets.Range("L16").Value = "instrument_id,session,date"
$datasets.Range("M16").Value = "instrument_id: ES`nthreshold_method: fixed_count`nmin_size: 100"

# The multi-line notes cells otherwise trigger Excel's auto row-height
# (like a real wrapped-text edit would); re-fit back to the sheet's
# normal single-line row height to match the rest of the sheet.
$datasets.Rows.Item(14).AutoFit()
$datasets.Rows.Item(16).AutoFit()

# ---------------------------------------------------------------------
# INSTRUMENTS sheet — insert 3 new columns in front of the "notes"
# column (P), shifting notes/volume_col/units right to S/T/U, then fill
# in the header labels and the new values for the ES row (row 26).
# ---------------------------------------------------------------------
$instruments = $wb.Worksheets.Item("INSTRUMENTS")

$instruments.Range("P1:R1").EntireColumn.Insert()
# The insert copies the bold header style from column O; the new
# columns should be unstyled like the rest of the (non-"notes") headers.
$instruments.Range("P1:R1").ClearFormats()

$instruments.Range("P1").Value = "big_trades_dataset_id"
$instruments.Range("Q1").Value = "big_trades_proxy_dataset_id"
$instruments.Range("R1").Value = "big_trades_source_mode"

$instruments.Range("P26").Value = "ES_BIG_TRADES"
$instruments.Range("Q26").Value = "ES_BIG_TRADES_PROXY"
$instruments.Range("R26").Value = "real_then_proxy"
